$d = $word.ActiveDocument

# 1. Merge "No external tools are used as like " + "Jira" (spell-checked run) + " or "
#    into a single run, dropping the spellStart/spellEnd proofErr markers that wrapped "Jira".
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("No external tools are used as like Jira or ", $true, $false, $false, $false, $false, $true, 1, $false, "No external tools are used as like Jira or ", 2)

# 2. Put the cursor right after "selenium " (before the comma) and drop the "_GoBack" bookmark
#    there -- this splits the trailing "," off into its own run, matching a fresh Word edit.
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("selenium ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng3)

# 3. Touch the table's "apply style" flags so Word materializes the full set of wdtblLook
#    boolean attributes (firstRow/lastRow/firstColumn/lastColumn/noHBand/noVBand) instead of
#    just the legacy w:val bitmask.
$tbl = $d.Tables(1)
$tbl.ApplyStyleHeadingRows = $false
$tbl.ApplyStyleLastRow = $false
$tbl.ApplyStyleFirstColumn = $false
$tbl.ApplyStyleLastColumn = $false
$tbl.ApplyStyleRowBands = $true
$tbl.ApplyStyleColumnBands = $true

